$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 657 (shifts existing rows 657..698 down to 658..699)
$ws.Rows.Item(657).Insert()

# Write the new row's data. The date column must stay plain text (not get
# auto-converted to a date serial by Excel's "smart" literal parsing), so we
# write it as a formula-computed string first and then convert that formula
# to a static value via copy / paste-values (mirrors how a user would type
# the value and then "paste as text" - it does not trigger date detection).
$ws.Cells.Item(657, 1).Formula = "=""2026/01/15"""
$ws.Cells.Item(657, 1).Copy()
$ws.Cells.Item(657, 1).PasteSpecial(-4163)

$ws.Cells.Item(657, 2).Value = "木"
$ws.Cells.Item(657, 3).Value = 20
$ws.Cells.Item(657, 4).Value = 201
